# soft delete added & phone number validation added
# - Update the phone number in A6 so it reflects an 11-digit value
#   (5555555555 -> 55555555555), matching the new phone number
#   validation rules.
# - Move the active selection to C12 (where the next phone number
#   entry / soft-delete flag would be added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're working on the active sheet before changing the
# selection, so the saved sheetView reflects the new active cell.
$ws.Activate()

# Phone number validation: the value in A6 grew by one digit.
$ws.Range("A6").Value = 55555555555

# Update the selected/active cell shown when the workbook is reopened.
$ws.Range("C12").Select()
